$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 2 (shifts existing data down by two rows)
$ws.Rows("2:3").Insert()

# Row 2: new "Posdoctorado" entry
$ws.Range("A2").Value = "Posdoctorado"
$ws.Range("B2").Value = "Desde 2023"
$ws.Range("C2").Value = "Asociación Red de Mujeres Víctimas y Profesionales"
$ws.Range("D2").Value = "Bogotá, Colombia"
$ws.Range("E2").Value = "\textbf{Proyecto: } La necesidad de generar procesos de reparación social a las mujeres víctimas y sobrevivientes de violencias sexuales en el marco del conflicto armado desde el quehacer periodístico. Diversas propuestas de tratamiento según contextos"

# Row 3: continuation detail for the new entry
$ws.Range("E3").Value = "Financiación del Ministerio de Ciencia Tecnología e Innovación - Minciencias"

# Append a new row 22 with another work experience entry
$ws.Range("A22").Value = "Comunicadora"
$ws.Range("B22").Value = "Mar 2006 - Nov 2009"
$ws.Range("C22").Value = "Brújula Comunicaciones "
$ws.Range("D22").Value = "Bogotá, Colombia"
$ws.Range("E22").Value = "Trabajo con medios de comunicación, realización de estrategias para acceder a derechos fundamentales"

$ws.Range("B22").Select()
